# Update odds values on row 4 of the active worksheet to reflect the
# latest FlashScore snapshot for the "Correcaminos - Atl. Morelia" match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G4"  = 2.77
    "H4"  = 3.1
    "J4"  = 3.25
    "L4"  = 3
    "P4"  = 2.95
    "Q4"  = 1.93
    "R4"  = 1.78
    "V4"  = 1.98
    "W4"  = 9.25
    "X4"  = 15
    "AA4" = 23
    "AC4" = 9.5
    "AD4" = 6.1
    "AE4" = 13
    "AH4" = 8
    "AI4" = 12
    "AL4" = 20
    "AM4" = 29
    "AN4" = 4.8
    "AO4" = 14.5
    "AP4" = 19.5
    "AQ4" = 60
    "AR4" = 80
    "AU4" = 6.4
    "AY4" = 19
    "BA4" = 80
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
